$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 177.8
$ws.Range("I5").Value = 133.33333
$ws.Range("J5").Value = 244.5
$ws.Range("K5").Value = 133.33333
$ws.Range("L5").Value = 244.5
$ws.Range("M5").Value = -18.33332999999999
$ws.Range("N5").Value = -474.5

$ws.Range("H40").Value = 999.9583
$ws.Range("I40").Value = 974.75
$ws.Range("J40").Value = 1005
$ws.Range("K40").Value = 974.75
$ws.Range("L40").Value = 1005
$ws.Range("M40").Value = -799.75
$ws.Range("N40").Value = -1355

$ws.Range("H74").Value = 3900
$ws.Range("I74").Value = 3836.3635
$ws.Range("K74").Value = 3836.3635
$ws.Range("M74").Value = -2900.3635

$ws.Range("H77").Value = 3900
$ws.Range("I77").Value = 3836.3635
$ws.Range("K77").Value = 19181.8175
$ws.Range("M77").Value = -14501.8175

$ws.Range("H129").Value = 832.4559
$ws.Range("J129").Value = 972.9804
$ws.Range("L129").Value = 2918.9412
$ws.Range("N129").Value = -12918.9412

$ws.Range("H138").Value = 1076295.2
$ws.Range("I138").Value = 1758.9524
$ws.Range("J138").Value = 1882197.5
$ws.Range("K138").Value = 5276.857199999999
$ws.Range("L138").Value = 5646592.5
$ws.Range("M138").Value = -136.8571999999995
$ws.Range("N138").Value = -5656872.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1838.5
$ws.Range("I2").Value = 1859.6666
$ws.Range("J2").Value = 1775
$ws.Range("K2").Value = 1859.6666
$ws.Range("L2").Value = 1775
$ws.Range("M2").Value = -1746.6666
$ws.Range("N2").Value = -2001

$ws.Range("H5").Value = 203.1
$ws.Range("J5").Value = 340
$ws.Range("L5").Value = 340
$ws.Range("N5").Value = -564

$ws.Range("H32").Value = 15339.419
$ws.Range("I32").Value = 15652.356
$ws.Range("J32").Value = 13582.154
$ws.Range("K32").Value = 15652.356
$ws.Range("L32").Value = 13582.154
$ws.Range("M32").Value = -15365.356
$ws.Range("N32").Value = -14156.154

$ws.Range("H61").Value = 166842780
$ws.Range("I61").Value = 125126420
$ws.Range("J61").Value = 250275500
$ws.Range("K61").Value = 125126420
$ws.Range("L61").Value = 250275500
$ws.Range("M61").Value = -125126208
$ws.Range("N61").Value = -250275924

$ws.Range("H74").Value = 8400862
$ws.Range("I74").Value = 9297184
$ws.Range("J74").Value = 333971.34
$ws.Range("K74").Value = 9297184
$ws.Range("L74").Value = 333971.34
$ws.Range("M74").Value = -9296310
$ws.Range("N74").Value = -335719.34

$ws.Range("H77").Value = 8400862
$ws.Range("I77").Value = 9297184
$ws.Range("J77").Value = 333971.34
$ws.Range("K77").Value = 46485920
$ws.Range("L77").Value = 1669856.7
$ws.Range("M77").Value = -46481552
$ws.Range("N77").Value = -1678592.7

$ws.Range("H97").Value = 1087.6154
$ws.Range("I97").Value = 1164.9
$ws.Range("J97").Value = 830
$ws.Range("K97").Value = 1164.9
$ws.Range("L97").Value = 830
$ws.Range("M97").Value = -668.9000000000001
$ws.Range("N97").Value = -1822

$ws.Range("H116").Value = 1838.5
$ws.Range("I116").Value = 1859.6666
$ws.Range("J116").Value = 1775
$ws.Range("K116").Value = 1859.6666
$ws.Range("L116").Value = 1775
$ws.Range("M116").Value = 434.3334
$ws.Range("N116").Value = -6363

$ws.Range("H136").Value = 166842780
$ws.Range("I136").Value = 125126420
$ws.Range("J136").Value = 250275500
$ws.Range("K136").Value = 375379260
$ws.Range("L136").Value = 750826500
$ws.Range("M136").Value = -375376710
$ws.Range("N136").Value = -750831600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1838.5
$ws.Range("I3").Value = 1859.6666
$ws.Range("J3").Value = 1775
$ws.Range("K3").Value = 1859.6666
$ws.Range("L3").Value = 1775
$ws.Range("M3").Value = -1745.6666
$ws.Range("N3").Value = -2003

$ws.Range("H4").Value = 203.1
$ws.Range("J4").Value = 340
$ws.Range("L4").Value = 340
$ws.Range("N4").Value = -570

$ws.Range("H20").Value = 1952.9
$ws.Range("I20").Value = 1158.6
$ws.Range("K20").Value = 1158.6
$ws.Range("M20").Value = -911.5999999999999

$ws.Range("H134").Value = 2928.0588
$ws.Range("I134").Value = 3333.0833
$ws.Range("J134").Value = 1956
$ws.Range("K134").Value = 9999.249899999999
$ws.Range("L134").Value = 5868
$ws.Range("M134").Value = -7464.249899999999
$ws.Range("N134").Value = -10938

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3165.121
$ws.Range("I31").Value = 1146.2593
$ws.Range("J31").Value = 12250
$ws.Range("K31").Value = 1146.2593
$ws.Range("L31").Value = 12250
$ws.Range("M31").Value = -851.2592999999999
$ws.Range("N31").Value = -12840

$ws.Range("H34").Value = 3165.121
$ws.Range("I34").Value = 1146.2593
$ws.Range("J34").Value = 12250
$ws.Range("K34").Value = 1146.2593
$ws.Range("L34").Value = 12250
$ws.Range("M34").Value = -944.2592999999999
$ws.Range("N34").Value = -12654

$ws.Range("H62").Value = 3025.2
$ws.Range("I62").Value = 2930.3333
$ws.Range("J62").Value = 3167.5
$ws.Range("K62").Value = 2930.3333
$ws.Range("L62").Value = 3167.5
$ws.Range("M62").Value = -2306.3333
$ws.Range("N62").Value = -4415.5

$ws.Range("H65").Value = 3025.2
$ws.Range("I65").Value = 2930.3333
$ws.Range("J65").Value = 3167.5
$ws.Range("K65").Value = 14651.6665
$ws.Range("L65").Value = 15837.5
$ws.Range("M65").Value = -11531.6665
$ws.Range("N65").Value = -22077.5

$ws.Range("H122").Value = 1560.1
$ws.Range("I122").Value = 1050.125
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 3150.375
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -700.375
$ws.Range("N122").Value = -15700

$ws.Range("H132").Value = 38007.5
$ws.Range("I132").Value = 2115.4285
$ws.Range("J132").Value = 145683.72
$ws.Range("K132").Value = 6346.2855
$ws.Range("L132").Value = 437051.16
$ws.Range("M132").Value = -3816.2855
$ws.Range("N132").Value = -442111.16

$ws.Range("H134").Value = 56411.75
$ws.Range("I134").Value = 2900.1428
$ws.Range("J134").Value = 181272.17
$ws.Range("K134").Value = 8700.428400000001
$ws.Range("L134").Value = 543816.51
$ws.Range("M134").Value = -6165.428400000001
$ws.Range("N134").Value = -548886.51

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 880.13635
$ws.Range("J5").Value = 1115.8334
$ws.Range("L5").Value = 3347.5002
$ws.Range("N5").Value = -3571.5002

$ws.Range("H107").Value = 1312.5834
$ws.Range("I107").Value = 935.9231
$ws.Range("J107").Value = 1757.7273
$ws.Range("K107").Value = 2807.7693
$ws.Range("L107").Value = 5273.1819
$ws.Range("M107").Value = -887.7692999999999
$ws.Range("N107").Value = -9113.1819

$ws.Range("H122").Value = 717.9524
$ws.Range("I122").Value = 332.75
$ws.Range("J122").Value = 1231.5555
$ws.Range("K122").Value = 2994.75
$ws.Range("L122").Value = 11083.9995
$ws.Range("M122").Value = -544.75
$ws.Range("N122").Value = -15983.9995

$ws.Range("H132").Value = 1525.3334
$ws.Range("I132").Value = 896.6667
$ws.Range("J132").Value = 1944.4445
$ws.Range("K132").Value = 8070.0003
$ws.Range("L132").Value = 17500.0005
$ws.Range("M132").Value = -5540.0003
$ws.Range("N132").Value = -22560.0005

$ws.Range("H135").Value = 880.13635
$ws.Range("J135").Value = 1115.8334
$ws.Range("L135").Value = 10042.5006
$ws.Range("N135").Value = -15112.5006

$ws.Range("H137").Value = 26465.791
$ws.Range("I137").Value = 944.4545000000001
$ws.Range("J137").Value = 48060.77
$ws.Range("K137").Value = 2833.3635
$ws.Range("L137").Value = 144182.31
$ws.Range("M137").Value = 2266.6365
$ws.Range("N137").Value = -154382.31

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 99645.2
$ws.Range("J139").Value = 99645.2
$ws.Range("L139").Value = 99645.2
$ws.Range("N139").Value = -109925.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 746.3077
$ws.Range("I22").Value = 575
$ws.Range("J22").Value = 822.44446
$ws.Range("K22").Value = 575
$ws.Range("L22").Value = 822.44446
$ws.Range("M22").Value = -280
$ws.Range("N22").Value = -1412.44446

$ws.Range("H27").Value = 746.3077
$ws.Range("I27").Value = 575
$ws.Range("J27").Value = 822.44446
$ws.Range("K27").Value = 575
$ws.Range("L27").Value = 822.44446
$ws.Range("M27").Value = -468
$ws.Range("N27").Value = -1036.44446

$ws.Range("H46").Value = 1024.8334
$ws.Range("I46").Value = 578.7143
$ws.Range("J46").Value = 1308.7273
$ws.Range("K46").Value = 578.7143
$ws.Range("L46").Value = 1308.7273
$ws.Range("M46").Value = -390.7143
$ws.Range("N46").Value = -1684.7273

$ws.Range("H132").Value = 79130.62
$ws.Range("I132").Value = 1188.6666
$ws.Range("K132").Value = 3565.9998
$ws.Range("M132").Value = -1035.9998

$ws.Range("H136").Value = 183074.55
$ws.Range("I136").Value = 334266.66
$ws.Range("J136").Value = 126377.5
$ws.Range("K136").Value = 1002799.98
$ws.Range("L136").Value = 379132.5
$ws.Range("M136").Value = -1000249.98
$ws.Range("N136").Value = -384232.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 155920.23
$ws.Range("I132").Value = 200500
$ws.Range("J132").Value = 128057.875
$ws.Range("K132").Value = 601500
$ws.Range("L132").Value = 384173.625
$ws.Range("M132").Value = -598970
$ws.Range("N132").Value = -389233.625

$ws.Range("H136").Value = 58332.113
$ws.Range("I136").Value = 34440.668
$ws.Range("J136").Value = 201680.8
$ws.Range("K136").Value = 103322.004
$ws.Range("L136").Value = 605042.3999999999
$ws.Range("M136").Value = -100772.004
$ws.Range("N136").Value = -610142.3999999999

$ws.Range("H137").Value = 42400.625
$ws.Range("J137").Value = 42400.625
$ws.Range("L137").Value = 42400.625
$ws.Range("N137").Value = -52600.625
